# #3550 fix typo in trx import sheet
# Rename the misspelled "Pament-Import" sheet to "Payment-Import" and
# make it the active/selected sheet (matching the author's workflow of
# fixing the typo while that sheet was in focus).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Pament-Import")
$ws.Name = "Payment-Import"

# Make the corrected sheet the active tab.
$ws.Activate()

# Restore the frozen-pane's active cell to A2 (the sheet's existing
# selection target) so only the tab-activation state changes.
$ws.Range("A2").Select()
